# Update cached leve-profit calculation values across multiple crafting-class sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H:N), per scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road (Potion)
$ws.Range("H17").Value = 3129461.5
$ws.Range("I17").Value = 1087.3334
$ws.Range("J17").Value = 3453086.5
$ws.Range("K17").Value = 3262.0002
$ws.Range("L17").Value = 10359259.5
$ws.Range("M17").Value = -3094.0002
$ws.Range("N17").Value = -10359595.5

# Row 40: Stuck in the Moment (Horn Glue)
$ws.Range("H40").Value = 2315.4
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2315.4
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2315.4
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2665.4

# Row 131: Mindful Study (Grade 5 Tincture of Mind)
$ws.Range("H131").Value = 1965
$ws.Range("I131").Value = 1965
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 5895
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -855

# Row 135: For Tired Minds (Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 12503531
$ws.Range("I135").Value = 638.4545000000001
$ws.Range("J135").Value = 71445736
$ws.Range("K135").Value = 5746.0905
$ws.Range("L135").Value = 643011624
$ws.Range("M135").Value = -3211.0905

# Row 138: All-night Crafting (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 11113487
$ws.Range("I138").Value = 35715520
$ws.Range("J138").Value = 2891.6936
$ws.Range("K138").Value = 107146560
$ws.Range("L138").Value = 8675.0808
$ws.Range("M138").Value = -107141420
$ws.Range("N138").Value = -18955.0808

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks (Mythril Ingot)
$ws.Range("H45").Value = 2330.5642
$ws.Range("I45").Value = 2114.8965
$ws.Range("J45").Value = 2956
$ws.Range("K45").Value = 2114.8965
$ws.Range("L45").Value = 2956
$ws.Range("M45").Value = -1737.8965
$ws.Range("N45").Value = -3710

# Row 122: Haste for High Durium (High Durium Nugget)
$ws.Range("H122").Value = 2689.2778
$ws.Range("I122").Value = 2107.3076
$ws.Range("J122").Value = 4202.4
$ws.Range("K122").Value = 6321.9228
$ws.Range("L122").Value = 12607.2
$ws.Range("M122").Value = -3871.9228

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt (Iron Ingot)
$ws.Range("H20").Value = 1225
$ws.Range("I20").Value = 1253.2
$ws.Range("J20").Value = 1178
$ws.Range("K20").Value = 1253.2
$ws.Range("L20").Value = 1178
$ws.Range("M20").Value = -1006.2
$ws.Range("N20").Value = -1672

# Row 22: Riveting Run (Iron Rivets)
$ws.Range("H22").Value = 272.80768
$ws.Range("I22").Value = 240.8
$ws.Range("J22").Value = 379.5
$ws.Range("K22").Value = 240.8
$ws.Range("L22").Value = 379.5
$ws.Range("M22").Value = -67.80000000000001
$ws.Range("N22").Value = -725.5

# Row 57: No Refunds, Only Exchanges (Cobalt File)
$ws.Range("H57").Value = 31591.666
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 31591.666
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 31591.666
$ws.Range("N57").Value = -33031.666

# Row 86: Through Thick and Thin (Adamantite Nugget)
$ws.Range("H86").Value = 1553.9512
$ws.Range("I86").Value = 1386.591
$ws.Range("J86").Value = 1747.7368
$ws.Range("K86").Value = 1386.591
$ws.Range("L86").Value = 1747.7368
$ws.Range("M86").Value = -263.5909999999999

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (Adamantite Nugget)
$ws.Range("H89").Value = 1553.9512
$ws.Range("I89").Value = 1386.591
$ws.Range("J89").Value = 1747.7368
$ws.Range("K89").Value = 6932.955
$ws.Range("L89").Value = 8738.683999999999
$ws.Range("M89").Value = -1316.955

# Row 94: High Steal (High Steel Nugget)
$ws.Range("H94").Value = 1235.3334
$ws.Range("I94").Value = 832.7692
$ws.Range("J94").Value = 2282
$ws.Range("K94").Value = 832.7692
$ws.Range("L94").Value = 2282
$ws.Range("M94").Value = -381.7692

# Row 105: Ingot to Wing It (Molybdenum Ingot)
$ws.Range("H105").Value = 2719.9412
$ws.Range("I105").Value = 2719.9167
$ws.Range("J105").Value = 2720
$ws.Range("K105").Value = 2719.9167
$ws.Range("L105").Value = 2720
$ws.Range("M105").Value = -972.9167000000002
$ws.Range("N105").Value = -6214

# Row 134: Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 25747.955
$ws.Range("I134").Value = 27339.268
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 82017.804
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -79482.804

# Row 136: Maintaining the Maintainers (Cobalt Tungsten File)
$ws.Range("H136").Value = 31591.666
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 31591.666
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 31591.666
$ws.Range("N136").Value = -41791.666

$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back (Square Maple Shield)
$ws.Range("H6").Value = 66667500
$ws.Range("I6").Value = 30000000
$ws.Range("J6").Value = 77143930
$ws.Range("K6").Value = 30000000
$ws.Range("L6").Value = 77143930
$ws.Range("M6").Value = -29999887

# Row 20: Re-crating the Scene (Iron Spear)
$ws.Range("H20").Value = 49000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 49000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 49000
$ws.Range("N20").Value = -49472

# Row 30: Polearms Aplenty (Iron Spear)
$ws.Range("H30").Value = 49000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 49000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 49000
$ws.Range("N30").Value = -49182

# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 3175.139
$ws.Range("I31").Value = 1734.6552
$ws.Range("J31").Value = 9142.857
$ws.Range("K31").Value = 1734.6552
$ws.Range("L31").Value = 9142.857
$ws.Range("M31").Value = -1439.6552
$ws.Range("N31").Value = -9732.857

# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 3175.139
$ws.Range("I34").Value = 1734.6552
$ws.Range("J34").Value = 9142.857
$ws.Range("K34").Value = 1734.6552
$ws.Range("L34").Value = 9142.857
$ws.Range("M34").Value = -1532.6552
$ws.Range("N34").Value = -9546.857

# Row 122: Timber of Tenkonto (Horse Chestnut Lumber)
$ws.Range("H122").Value = 3583.5
$ws.Range("I122").Value = 4375.25
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 13125.75
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -10675.75

# Row 128: An A-prop-riate Request (Ironwood Spear)
$ws.Range("H128").Value = 49000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49000
$ws.Range("N128").Value = -58960

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch (Chamomile Tea)
$ws.Range("H34").Value = 1689
$ws.Range("I34").Value = 597.3333
$ws.Range("J34").Value = 2780.6667
$ws.Range("K34").Value = 1791.9999
$ws.Range("L34").Value = 8342.000100000001
$ws.Range("M34").Value = -1707.9999
$ws.Range("N34").Value = -8510.000100000001

# Row 37: I Love Lamprey (Eel Pie)
$ws.Range("H37").Value = 62571000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 62571000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 187713000
$ws.Range("N37").Value = -187713224

# Row 63: The Next to Last Supper (Stuffed Cabbage Rolls)
$ws.Range("H63").Value = 5078.5
$ws.Range("I63").Value = 2900
$ws.Range("J63").Value = 7257
$ws.Range("K63").Value = 8700
$ws.Range("L63").Value = 21771
$ws.Range("M63").Value = -7951
$ws.Range("N63").Value = -23269

# Row 64: The Aroma of Faith (Baked Onion Soup)
$ws.Range("H64").Value = 3604
$ws.Range("I64").Value = 2906
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 8718
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -8448
$ws.Range("N64").Value = -15540

# Row 66: Nostalgia through the Stomach (L) (Stuffed Cabbage Rolls)
$ws.Range("H66").Value = 5078.5
$ws.Range("I66").Value = 2900
$ws.Range("J66").Value = 7257
$ws.Range("K66").Value = 26100
$ws.Range("L66").Value = 65313
$ws.Range("M66").Value = -22356
$ws.Range("N66").Value = -72801

# Row 67: Soup's On (L) (Baked Onion Soup)
$ws.Range("H67").Value = 3604
$ws.Range("I67").Value = 2906
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 8718
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -7782
$ws.Range("N67").Value = -16872

# Row 106: Herky Jerky (Jerked Jhammel)
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 120: A Happy End (Paella)
$ws.Range("H120").Value = 17151.428
$ws.Range("I120").Value = 10030
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 30090
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -25252

# Row 121: A Cookie for Your Troubles (Coffee Biscuit)
$ws.Range("H121").Value = 3472.1143
$ws.Range("I121").Value = 440.625
$ws.Range("J121").Value = 4370.3335
$ws.Range("K121").Value = 1321.875
$ws.Range("L121").Value = 13111.0005
$ws.Range("M121").Value = -11.875
$ws.Range("N121").Value = -15731.0005

# Row 129: Comfort Food (Yakow Moussaka)
$ws.Range("H129").Value = 209862.92
$ws.Range("I129").Value = 881.5
$ws.Range("J129").Value = 418844.34
$ws.Range("K129").Value = 2644.5
$ws.Range("L129").Value = 1256533.02
$ws.Range("M129").Value = 2355.5
$ws.Range("N129").Value = -1266533.02

# Row 136: Simple Is Hardest (Spaghetti al Olio e Peperoncino)
$ws.Range("H136").Value = 1858.6666
$ws.Range("I136").Value = 1076.4166
$ws.Range("J136").Value = 4987.6665
$ws.Range("K136").Value = 3229.2498
$ws.Range("L136").Value = 14962.9995
$ws.Range("M136").Value = 1870.7502
$ws.Range("N136").Value = -25162.9995

$ws = $wb.Worksheets.Item("GSM")
# Row 45: Lode It Up (Mythril Ring)
$ws.Range("H45").Value = 20217.334
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 20217.334
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 20217.334
$ws.Range("N45").Value = -21335.334

# Row 113: Copious Crystal Cannons (Manasilver Nugget)
$ws.Range("H113").Value = 4128.5713
$ws.Range("I113").Value = 3100
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 3100
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -930
$ws.Range("N113").Value = -9840

# Row 132: On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 75355.81
$ws.Range("I132").Value = 56924.844
$ws.Range("J132").Value = 250450
$ws.Range("K132").Value = 170774.532
$ws.Range("L132").Value = 751350
$ws.Range("M132").Value = -168244.532
$ws.Range("N132").Value = -756410

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly (Raptor Leather)
$ws.Range("H61").Value = 6399.909
$ws.Range("I61").Value = 3366.6667
$ws.Range("J61").Value = 7537.375
$ws.Range("K61").Value = 3366.6667
$ws.Range("L61").Value = 7537.375
$ws.Range("M61").Value = -3164.6667

# Row 100: Tiger in the Sack (Tiger Leather)
$ws.Range("I100").Value = 1580.6
$ws.Range("J100").Value = 2218
$ws.Range("K100").Value = 1580.6
$ws.Range("L100").Value = 2218
$ws.Range("M100").Value = -1039.6
$ws.Range("N100").Value = -3300

# Row 113: Peace in Rest (Atrociraptor Leather)
$ws.Range("H113").Value = 6399.909
$ws.Range("I113").Value = 3366.6667
$ws.Range("J113").Value = 7537.375
$ws.Range("K113").Value = 3366.6667
$ws.Range("L113").Value = 7537.375
$ws.Range("M113").Value = -1196.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 27: Hitting Below the Belt (Cotton Breeches of Crafting)
$ws.Range("H27").Value = 33417.332
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 33417.332
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 33417.332
$ws.Range("N27").Value = -33555.332

# Row 113: A Tender Table (Pixie Floss)
$ws.Range("H113").Value = 2079832.2
$ws.Range("I113").Value = 1153.7778
$ws.Range("J113").Value = 6756859
$ws.Range("K113").Value = 3461.3334
$ws.Range("L113").Value = 20270577
$ws.Range("M113").Value = -1291.3334

# Row 115: Gloves Come in Handy (Pixie Cotton Sleeves of Crafting)
$ws.Range("H115").Value = 30377
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 30377
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 30377
$ws.Range("N115").Value = -33511

# Row 123: Helping Handwear (Fingerless Darkhempen Gloves of Healing)
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 126: A Polished Purchase (Snow Linen)
$ws.Range("H126").Value = 1118.84
$ws.Range("I126").Value = 1068.4546
$ws.Range("J126").Value = 1488.3334
$ws.Range("K126").Value = 3205.3638
$ws.Range("L126").Value = 4465.0002
$ws.Range("M126").Value = -735.3638000000001
$ws.Range("N126").Value = -9405.0002

# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 20647424
$ws.Range("I136").Value = 25177992
$ws.Range("J136").Value = 8167.222
$ws.Range("K136").Value = 75533976
$ws.Range("L136").Value = 24501.666
$ws.Range("M136").Value = -75531426
$ws.Range("N136").Value = -29601.666
